$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto market data values (prices & 1h volume change %).
# Cells whose new text would otherwise be auto-parsed as a plain number by Excel
# get NumberFormat "@" (Text) first, so they are stored as text, matching the source data.

# Row 2
$ws.Range("D2").Value = "42.483.89"
$ws.Range("E2").Value = "  -6.82%  "

# Row 3
$ws.Range("D3").Value = "2.204.47"
$ws.Range("E3").Value = "  -7.41%  "

# Row 4
$ws.Range("E4").Value = "  -0.21%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.62"
$ws.Range("E5").Value = "  -1.62%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.74"
$ws.Range("E6").Value = "  -12.46%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.572"
$ws.Range("E7").Value = "  -9.36%  "

# Row 8
$ws.Range("E8").Value = "  +0.09%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.552"
$ws.Range("E9").Value = "  -10.83%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.41"
$ws.Range("E10").Value = "  -11.65%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.62"
$ws.Range("E11").Value = "  -4.00%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0821"
$ws.Range("E12").Value = "  -11.22%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.53"
$ws.Range("E13").Value = "  -12.27%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.105"
$ws.Range("E14").Value = "  -4.49%  "

# Row 15
$ws.Range("D15").Value = "2.541.79"
$ws.Range("E15").Value = "  -7.42%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.850"
$ws.Range("E16").Value = "  -14.23%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.82"
$ws.Range("E17").Value = "  -11.41%  "

# Row 18
$ws.Range("D18").Value = "2.186.26"
$ws.Range("E18").Value = "  -8.15%  "

# Row 19
$ws.Range("D19").Value = "42.360.14"
$ws.Range("E19").Value = "  -6.97%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.17"
$ws.Range("E20").Value = "  +7.57%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.47"
$ws.Range("E21").Value = "  -12.05%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0941"
$ws.Range("E22").Value = "  -12.24%  "

# Row 23
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.61"
$ws.Range("E23").Value = "  -12.77%  "

# Row 24
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.16"
$ws.Range("E24").Value = "  -9.56%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "231.87"
$ws.Range("E25").Value = "  -11.83%  "

# Row 26
$ws.Range("E26").Value = "  -10.11%  "

# Row 27
$ws.Range("E27").Value = "  +0.20%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.04"
$ws.Range("E28").Value = "  -10.24%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.14"
$ws.Range("E29").Value = "  -9.02%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.43"
$ws.Range("E30").Value = "  -14.54%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.29"
$ws.Range("E31").Value = "  -10.27%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "157.37"
$ws.Range("E32").Value = "  -7.77%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0859"
$ws.Range("E33").Value = "  -11.31%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "33.35"
$ws.Range("E34").Value = "  -12.76%  "

# Row 35
$ws.Range("E35").Value = "  -9.39%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.20"
$ws.Range("E36").Value = "  +6.82%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.121"
$ws.Range("E37").Value = "  -7.81%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.43"
$ws.Range("E38").Value = "  -8.66%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.83"
$ws.Range("E39").Value = "  +5.13%  "

# Row 40
$ws.Range("E40").Value = "  -12.22%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.50"
$ws.Range("E41").Value = "  -12.83%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0315"
$ws.Range("E42").Value = "  -11.85%  "

# Row 43
$ws.Range("E43").Value = "  +0.02%  "

# Row 44
$ws.Range("D44").Value = "1.769.21"
$ws.Range("E44").Value = "  +7.26%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "88.49"
$ws.Range("E45").Value = "  -13.58%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.88"
$ws.Range("E46").Value = "  -10.95%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.203"
$ws.Range("E47").Value = "  -14.27%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "76.61"
$ws.Range("E48").Value = "  -8.64%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.26"
$ws.Range("E49").Value = "  -5.50%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "59.58"
$ws.Range("E50").Value = "  -15.31%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.43"
$ws.Range("E51").Value = "  -10.10%  "
